# "cap nhat msv Phan Dang Hung" - fill in the previously-empty MSV (student
# id) cell for row 3 ("Phan Dang Hung") of the team table with "2051063664".

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 4 of the table = "3 | Phan Đăng Hùng | (MSV) | Lên ý tưởng";
# column 3 is the MSV cell, which is currently empty (just the paragraph
# mark) for this row.
$cell = $t.Cell(4, 3)
$rng = $cell.Range
$rng.SetRange($rng.Start, $rng.End - 1)
$rng.InsertBefore("2051063664")

# Match the formatting used by the sibling MSV cells in this table
# (Times New Roman, 15pt / sz=30) by running a formatted find/replace over
# the text we just inserted.
$find = $rng.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Name = "Times New Roman"
$find.Replacement.Font.NameAscii = "Times New Roman"
$find.Replacement.Font.NameFarEast = "Times New Roman"
$find.Replacement.Font.NameBi = "Times New Roman"
$find.Replacement.Font.NameOther = "Times New Roman"
$find.Replacement.Font.Size = 15
$null = $find.Execute("2051063664", $false, $false, $false, $false, $false, $true, 1, $false, "2051063664", 2)
